# Applies the "Add data for 2024-10-28" update to the violent-crime-full-year
# workbook. This populates a new day's worth of incident counts into the
# 2024 (column K) totals across the Citywide Totals, By Neighborhood, and
# per-neighborhood sheets (and a couple of minor historical revisions to
# 2015 figures, column B, that came bundled with the same data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6701
$ws.Range("K3").Value = 6916
$ws.Range("B4").Value = 1705
$ws.Range("K4").Value = 1434
$ws.Range("K5").Value = 500
$ws.Range("K6").Value = 7593
$ws.Range("B7").Value = 23338
$ws.Range("K7").Value = 23144

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K4").Value = 15
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 297

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 420
$ws.Range("K3").Value = 462
$ws.Range("K6").Value = 502
$ws.Range("K7").Value = 1517

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 500

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 130
$ws.Range("K3").Value = 133
$ws.Range("K7").Value = 381

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 225
$ws.Range("K6").Value = 234
$ws.Range("K7").Value = 785

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 197
$ws.Range("K7").Value = 540

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 103
$ws.Range("K7").Value = 392

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 202
$ws.Range("K7").Value = 699
$ws.Range("K8").Value = 1517
$ws.Range("K9").Value = 106
$ws.Range("K11").Value = 429
$ws.Range("K15").Value = 244
$ws.Range("K19").Value = 678
$ws.Range("K20").Value = 558
$ws.Range("K25").Value = 109
$ws.Range("K27").Value = 217
$ws.Range("K29").Value = 1256
$ws.Range("K36").Value = 294
$ws.Range("K37").Value = 785
$ws.Range("K42").Value = 854
$ws.Range("K44").Value = 193
$ws.Range("K48").Value = 294
$ws.Range("K49").Value = 125
$ws.Range("K52").Value = 613
$ws.Range("K53").Value = 297
$ws.Range("B63").Value = 409
$ws.Range("K63").Value = 61
$ws.Range("K64").Value = 144
$ws.Range("K65").Value = 540
$ws.Range("K66").Value = 70
$ws.Range("K67").Value = 904
$ws.Range("K71").Value = 72
$ws.Range("K78").Value = 264
$ws.Range("K79").Value = 581
$ws.Range("K83").Value = 500
$ws.Range("K85").Value = 1066
$ws.Range("K86").Value = 142
$ws.Range("K88").Value = 249
$ws.Range("K91").Value = 277
$ws.Range("K95").Value = 381
$ws.Range("K97").Value = 181
$ws.Range("K98").Value = 117
$ws.Range("K99").Value = 392
$ws.Range("B101").Value = 23338
$ws.Range("K101").Value = 23144

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 246
$ws.Range("K3").Value = 330
$ws.Range("K6").Value = 258
$ws.Range("K7").Value = 904

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 446
$ws.Range("K6").Value = 366
$ws.Range("K7").Value = 1256

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 294

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 203
$ws.Range("K5").Value = 20
$ws.Range("K6").Value = 225
$ws.Range("K7").Value = 678

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 230
$ws.Range("K6").Value = 317
$ws.Range("K7").Value = 854

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 79
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 70
$ws.Range("K3").Value = 131
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 189
$ws.Range("K5").Value = 19
$ws.Range("K7").Value = 581

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 178
$ws.Range("K7").Value = 558

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 90
$ws.Range("K7").Value = 294

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 229
$ws.Range("K7").Value = 699

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 91
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 244

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 150
$ws.Range("K3").Value = 111
$ws.Range("K7").Value = 429

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 76
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K6").Value = 259
$ws.Range("K7").Value = 1066

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 173
$ws.Range("K7").Value = 613
